$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range('D2').Value = '49.856.39'
$ws.Range('E2').Value = '  +3.73%  '
$ws.Range('D3').Value = '2.647.03'
$ws.Range('E3').Value = '  +6.02%  '
Set-TextValue $ws 'D4' '1.00'
$ws.Range('E4').Value = '  +0.08%  '
Set-TextValue $ws 'D5' '114.04'
$ws.Range('E5').Value = '  +7.84%  '
Set-TextValue $ws 'D6' '326.10'
$ws.Range('E6').Value = '  +2.04%  '
$ws.Range('E7').Value = '  +1.50%  '
Set-TextValue $ws 'D8' '1.00'
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('E9').Value = '  +2.93%  '
Set-TextValue $ws 'D10' '41.04'
$ws.Range('E10').Value = '  +6.04%  '
Set-TextValue $ws 'D11' '20.11'
$ws.Range('E11').Value = '  +0.45%  '
$ws.Range('E13').Value = '  +1.01%  '
Set-TextValue $ws 'D14' '7.36'
$ws.Range('E14').Value = '  +3.74%  '
$ws.Range('D15').Value = '3.062.94'
$ws.Range('E15').Value = '  +6.04%  '
$ws.Range('D16').Value = '2.648.59'
$ws.Range('E16').Value = '  +6.07%  '
$ws.Range('E17').Value = '  +4.40%  '
$ws.Range('D18').Value = '49.781.47'
$ws.Range('E18').Value = '  +3.95%  '
Set-TextValue $ws 'D19' '13.16'
$ws.Range('E19').Value = '  +0.66%  '
$ws.Range('B20').Value = 'ImmutableX'
$ws.Range('C20').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws 'D20' '2.95'
$ws.Range('E20').Value = '  -0.05%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue $ws 'D21' '6.78'
$ws.Range('E21').Value = '  +2.07%  '
$ws.Range('D22').Value = '0.0₃0956'
$ws.Range('E22').Value = '  +2.58%  '
Set-TextValue $ws 'D23' '72.04'
$ws.Range('E23').Value = '  +1.24%  '
Set-TextValue $ws 'D24' '277.28'
$ws.Range('E24').Value = '  +1.94%  '
Set-TextValue $ws 'D25' '2.58'
$ws.Range('E25').Value = '  +2.62%  '
Set-TextValue $ws 'D26' '26.77'
$ws.Range('E26').Value = '  +3.96%  '
$ws.Range('E28').Value = '  +2.69%  '
$ws.Range('E29').Value = '  -1.22%  '
Set-TextValue $ws 'D30' '36.06'
$ws.Range('E30').Value = '  +3.88%  '
Set-TextValue $ws 'D31' '0.140'
$ws.Range('E31').Value = '  +0.21%  '
Set-TextValue $ws 'D32' '50.33'
$ws.Range('E32').Value = '  +2.55%  '
Set-TextValue $ws 'D33' '5.42'
$ws.Range('E33').Value = '  +2.50%  '
Set-TextValue $ws 'D34' '19.53'
$ws.Range('E34').Value = '  +2.25%  '
$ws.Range('E35').Value = '  +4.18%  '
$ws.Range('E36').Value = '  +0.01%  '
Set-TextValue $ws 'D37' '2.07'
$ws.Range('E37').Value = '  +7.16%  '
Set-TextValue $ws 'D38' '4.86'
$ws.Range('E38').Value = '  +6.96%  '
Set-TextValue $ws 'D39' '3.08'
$ws.Range('E39').Value = '  +8.18%  '
Set-TextValue $ws 'D40' '126.59'
$ws.Range('E40').Value = '  +3.75%  '
$ws.Range('E41').Value = '  +1.59%  '
$ws.Range('E42').Value = '  +1.54%  '
Set-TextValue $ws 'D43' '22.00'
$ws.Range('E43').Value = '  -1.39%  '
Set-TextValue $ws 'D44' '0.0314'
$ws.Range('E44').Value = '  +3.40%  '
$ws.Range('D45').Value = '2.078.37'
$ws.Range('E45').Value = '  +3.99%  '
$ws.Range('E46').Value = '  +5.27%  '
$ws.Range('E48').Value = '  +4.61%  '
$ws.Range('E49').Value = '  +2.42%  '
Set-TextValue $ws 'D50' '5.37'
$ws.Range('E50').Value = '  +3.84%  '
Set-TextValue $ws 'D51' '60.32'
$ws.Range('E51').Value = '  +7.34%  '
